# Weekly update to the Acelga (Terminal Hortofrutícola Agro Chillán) sheet.
# A new week of data (2 rows: Primera / Segunda) is inserted right after the
# existing row 617, pushing all the following rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 618-619; everything currently at/after row 618
# (including the trailing rows that fall off the old A1:R649 range) shifts
# down to 620-651.
$ws.Range("A618:A619").EntireRow.Insert()

# New row 618 - "Primera" quality, week of 2023-12-18 (serial 45267)
$ws.Cells.Item(618, 1).Value = 7
$ws.Cells.Item(618, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(618, 3).Value = "Ñuble"
$ws.Cells.Item(618, 4).Value = 45267
$ws.Cells.Item(618, 5).Value = 16
$ws.Cells.Item(618, 6).Value = 100112009
$ws.Cells.Item(618, 7).Value = "Acelga"
$ws.Cells.Item(618, 8).Value = "Sin especificar"
$ws.Cells.Item(618, 9).Value = "Primera"
$ws.Cells.Item(618, 10).Value = 200
$ws.Cells.Item(618, 11).Value = 700
$ws.Cells.Item(618, 12).Value = 700
$ws.Cells.Item(618, 13).Value = 700
$ws.Cells.Item(618, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(618, 15).Value = "Región de Ñuble"
$ws.Cells.Item(618, 16).Value = 700
$ws.Cells.Item(618, 17).Value = 1
$ws.Cells.Item(618, 18).Value = "Hortaliza"

# New row 619 - "Segunda" quality, same week
$ws.Cells.Item(619, 1).Value = 7
$ws.Cells.Item(619, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(619, 3).Value = "Ñuble"
$ws.Cells.Item(619, 4).Value = 45267
$ws.Cells.Item(619, 5).Value = 16
$ws.Cells.Item(619, 6).Value = 100112009
$ws.Cells.Item(619, 7).Value = "Acelga"
$ws.Cells.Item(619, 8).Value = "Sin especificar"
$ws.Cells.Item(619, 9).Value = "Segunda"
$ws.Cells.Item(619, 10).Value = 300
$ws.Cells.Item(619, 11).Value = 500
$ws.Cells.Item(619, 12).Value = 500
$ws.Cells.Item(619, 13).Value = 500
$ws.Cells.Item(619, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(619, 15).Value = "Región de Ñuble"
$ws.Cells.Item(619, 16).Value = 500
$ws.Cells.Item(619, 17).Value = 1
$ws.Cells.Item(619, 18).Value = "Hortaliza"
